# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) marking Control (0) vs MDD (1) rows, and
# refreshes the refitted Prediction/Error/Cross-Entropy-Loss values for the
# 100-iteration block (rows 2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell H1, styled like the other header cells (G1) ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# --- Updated numeric results for the 100-iteration block (rows 2-11) ---
$ws.Range("D2").Value = 0.5949309237310512
$ws.Range("E2").Value = 0.5949309237310512

$ws.Range("D3").Value = 0.5980592011429327
$ws.Range("E3").Value = 0.5980592011429327

$ws.Range("D4").Value = 0.5082293098340988
$ws.Range("E4").Value = 0.5082293098340988

$ws.Range("D6").Value = 0.5270036130372806
$ws.Range("E6").Value = 0.5270036130372806

$ws.Range("D7").Value = 0.4975607586296235
$ws.Range("E7").Value = 0.5024392413703764

$ws.Range("D8").Value = 0.5166199789900719
$ws.Range("E8").Value = 0.4833800210099281

$ws.Range("D9").Value = 0.48275463887822
$ws.Range("E9").Value = 0.51724536112178

$ws.Range("D10").Value = 0.6370068013320941
$ws.Range("E10").Value = 0.3629931986679059

$ws.Range("D11").Value = 0.4932044412360124
$ws.Range("E11").Value = 0.5067955587639876
$ws.Range("F11").Value = 0.7717534303665161

# --- New "Label" column values: 0 = Control, 1 = MDD ---
$labels = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1;
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1;
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
